$wb = $excel.ActiveWorkbook

# The "建物" (building) sheet incorrectly had property_category = "land"
# for all of its rows. Fix it to "building" (issue #5).
$ws = $wb.Worksheets.Item("建物")

$ws.Range("I2:I5").Value = "building"
